$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'27.138.60"
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'1.678.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.02%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "'214.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.66%  "

# Row 6
$ws.Range("E6").Value = "  -0.13%  "

# Row 7
$ws.Range("E7").Value = "  +0.09%  "

# Row 8
$ws.Range("D8").Value = "'22.78"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.44%  "

# Row 9
$ws.Range("E9").Value = "  +2.29%  "

# Row 10
$ws.Range("E10").Value = "  -0.40%  "

# Row 11
$ws.Range("D11").Value = "'0.0891"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.27%  "

# Row 12
$ws.Range("D12").Value = "'1.914.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.06%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.22%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "'1.654.89"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.26%  "

# Row 15
$ws.Range("D15").Value = "'0.554"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.62%  "

# Row 16
$ws.Range("D16").Value = "'66.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.08%  "

# Row 17
$ws.Range("D17").Value = "'27.109.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.33%  "

# Row 18
$ws.Range("D18").Value = "'235.13"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.25%  "

# Row 19
$ws.Range("D19").Value = "'7.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.19%  "

# Row 20
$ws.Range("E20").Value = "  +0.30%  "

# Row 21
$ws.Range("E21").Value = "  +0.12%  "

# Row 22
$ws.Range("D22").Value = "'4.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.61%  "

# Row 23
$ws.Range("E23").Value = "  +2.76%  "

# Row 24
$ws.Range("D24").Value = "'2.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.40%  "

# Row 25
$ws.Range("D25").Value = "'148.55"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.44%  "

# Row 26
$ws.Range("E26").Value = "  +2.33%  "

# Row 27
$ws.Range("E27").Value = "  -0.46%  "

# Row 28
$ws.Range("D28").Value = "'0.113"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.08%  "

# Row 29
$ws.Range("E29").Value = "  -0.02%  "

# Row 30
$ws.Range("E30").Value = "  +0.58%  "

# Row 31
$ws.Range("E31").Value = "  -0.46%  "

# Row 32
$ws.Range("E32").Value = "  -0.21%  "

# Row 33
$ws.Range("D33").Value = "'1.540.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.12%  "

# Row 34
$ws.Range("E34").Value = "  +1.17%  "

# Row 35
$ws.Range("E35").Value = "  -3.62%  "

# Row 36
$ws.Range("D36").Value = "'0.606"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.73%  "

# Row 37
$ws.Range("D37").Value = "'0.939"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.08%  "

# Row 38
$ws.Range("E38").Value = "  -0.13%  "

# Row 39
$ws.Range("E39").Value = "  -1.03%  "

# Row 40
$ws.Range("E40").Value = "  +2.20%  "

# Row 41
$ws.Range("E41").Value = "  +3.25%  "

# Row 42
$ws.Range("D42").Value = "'69.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.19%  "

# Row 43
$ws.Range("E43").Value = "  +0.11%  "

# Row 44
$ws.Range("E44").Value = "  -0.43%  "

# Row 45
$ws.Range("D45").Value = "'1.822.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.22%  "

# Row 46
$ws.Range("D46").Value = "'0.779"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.10%  "

# Row 47
$ws.Range("D47").Value = "'89.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.44%  "

# Row 48
$ws.Range("E48").Value = "  +6.03%  "

# Row 49
$ws.Range("E49").Value = "  +3.25%  "

# Row 50
$ws.Range("D50").Value = "'8.20"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.04%  "

# Row 51
$ws.Range("E51").Value = "  -0.26%  "
